# PlayerPerformance_3996.xlsx — additional scraping update
#
# 1. Add a new "Player Info" sheet (becomes the first sheet) with the
#    player's basic bio info.
# 2. On "ODI Batting" / "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE
#    and replace the full scorecard URL values with just the numeric
#    match code. Also drop the leftover blank INNING_NUMBER cells on
#    "ODI Batting" for innings the player didn't bat in.
# 3. Add a new "ODI Batting Extra" sheet (last sheet) with additional
#    per-match batting detail.

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------
# 1. "Player Info" sheet — inserted before "ODI Batting" so the final
#    sheet order is Player Info, ODI Batting, ODI Bowling, ...
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Reuse the bold/bordered header style already used by the other sheets.
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "3996"
$playerInfo.Range("B2").Value = "Kane W Richardson"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

# ---------------------------------------------------------------------
# 2a. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (header + values),
#     and remove the stray blank INNING_NUMBER cells.
# ---------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingRows = $battingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $battingRows; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Value
    if ($link) {
        $code = $link.ToString().Split("=")[-1]
        $cell.Value = $code
    }

    $inning = $battingSheet.Cells.Item($r, 2)
    if ([string]::IsNullOrEmpty($inning.Value)) {
        $inning.Value = ""
    }
}

# ---------------------------------------------------------------------
# 2b. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (header + values).
# ---------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingRows = $bowlingSheet.UsedRange.Rows.Count
for ($r = 2; $r -le $bowlingRows; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Value
    if ($link) {
        $code = $link.ToString().Split("=")[-1]
        $cell.Value = $code
    }
}

# ---------------------------------------------------------------------
# 3. "ODI Batting Extra" — new sheet appended at the end.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraSheet = $wb.Worksheets.Add($null, $lastSheet)
$extraSheet.Name = "ODI Batting Extra"

$extraSheet.Range("A1").Value = "MATCH_CODE"
$extraSheet.Range("B1").Value = "BATTING_POSITION"
$extraSheet.Range("C1").Value = "NUM_4"
$extraSheet.Range("D1").Value = "NUM_6"
$extraSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extraSheet.Range("F1").Value = "MAN_OF_MATCH"

$battingSheet.Range("A1").Copy()
$extraSheet.Range("A1:F1").PasteSpecial(-4122)

$extraRows = @(
    @("3679", 9,    $null, $null, $null,    "NO"),
    @("3681", 10,   "1",   "0",   "3.90%",  "NO"),
    @("3697", 9,    $null, $null, $null,    "NO"),
    @("3875", $null,$null, $null, $null,    "NO"),
    @("3876", 10,   $null, $null, $null,    "NO"),
    @("3877", $null,$null, $null, $null,    "NO"),
    @("3884", $null,$null, $null, $null,    "NO"),
    @("4069", 11,   $null, $null, $null,    "NO"),
    @("4071", $null,$null, $null, $null,    "NO"),
    @("4074", 10,   "0",   "0",   "0.47%",  "NO"),
    @("4166", $null,$null, $null, $null,    "NO"),
    @("4167", 9,    "0",   "1",   "6.83%",  "NO"),
    @("4170", $null,$null, $null, $null,    "NO"),
    @("4276", 9,    "1",   "0",   "1.53%",  "NO"),
    @("4277", 11,   "0",   "0",   "0.33%",  "NO"),
    @("4319", 11,   $null, $null, $null,    "NO"),
    @("4322", 10,   $null, $null, $null,    "NO"),
    @("4398", 10,   "4",   "1",   "7.89%",  "NO"),
    @("4399", $null,$null, $null, $null,    "NO"),
    @("4421", $null,$null, $null, $null,    $null)
)

$extraSheet.Range("C2:E21").NumberFormat = "@"

$r = 2
foreach ($row in $extraRows) {
    $extraSheet.Cells.Item($r, 1).Value = $row[0]

    if ($row[1] -ne $null) {
        $extraSheet.Cells.Item($r, 2).Value = $row[1]
    } else {
        $extraSheet.Cells.Item($r, 2).Value = ""
    }

    for ($c = 3; $c -le 5; $c++) {
        $v = $row[$c - 1]
        if ($v -ne $null) {
            $extraSheet.Cells.Item($r, $c).Value = $v
        } else {
            $extraSheet.Cells.Item($r, $c).Value = ""
        }
    }

    if ($row[5] -ne $null) {
        $extraSheet.Cells.Item($r, 6).Value = $row[5]
    } else {
        $extraSheet.Cells.Item($r, 6).Value = ""
    }

    $r++
}
